$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting existing rows 124:249 down to 125:250
$ws.Rows("124:124").Insert()

# Populate the new row 124 with the latest weekly price report
$ws.Range("A124").Value = 5
$ws.Range("B124").Value = "Macroferia Regional de Talca"
$ws.Range("C124").Value = "Maule"
$ws.Range("D124").Value = 44705
$ws.Range("E124").Value = 7
$ws.Range("F124").Value = 100112009
$ws.Range("G124").Value = "Acelga"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 500
$ws.Range("K124").Value = 3000
$ws.Range("L124").Value = 3000
$ws.Range("M124").Value = 3000
$ws.Range("N124").Value = "$/docena de atados (4 kilos)"
$ws.Range("O124").Value = "Región del Maule"
$ws.Range("P124").Value = 750
$ws.Range("Q124").Value = 4
$ws.Range("R124").Value = "Hortaliza"
